$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 5 (columns D:G) with "ok", matching the existing C5 value,
# so the whole row C5:G5 reads "ok".
$ws.Range("D5").Value = "ok"
$ws.Range("E5").Value = "ok"
$ws.Range("F5").Value = "ok"
$ws.Range("G5").Value = "ok"

# Move the active selection from D5 to G5.
$ws.Range("G5").Select()
